$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LargeBin")

# Fix I15 to use the same style as H15 (removing a duplicate/unused style)
$ws.Range("H15").Copy()
$ws.Range("I15").PasteSpecial(-4122)

# Fill in the previously-empty I/J values for rows 17-44
$values = @{
    17 = @(0.78200000000000003, 0.67500000000000004)
    18 = @(0.99399999999999999, 0.025000000000000001)
    19 = @(0.88600000000000001, 0.35399999999999998)
    20 = @(0.77300000000000002, 0.33900000000000002)
    21 = @(0.89100000000000001, 0.30399999999999999)
    22 = @(0.89400000000000002, 0.159)
    23 = @(0.78500000000000003, 0.44800000000000001)
    24 = @(0.1, 0)
    25 = @(0.85199999999999998, 0.47199999999999998)
    26 = @(1, 0)
    27 = @(0.75900000000000001, 0.50900000000000001)
    28 = @(1, 0)
    29 = @(0.877, 0.41299999999999998)
    30 = @(1, 0)
    31 = @(0.78400000000000003, 0.501)
    32 = @(1, 0)
    33 = @(1, 0.0030000000000000001)
    34 = @(0.89600000000000002, 0.153)
    35 = @(1, 0.0070000000000000001)
    36 = @(0.88400000000000001, 0.185)
    37 = @(0.79400000000000004, 0.70399999999999996)
    38 = @(1, 0)
    39 = @(1, 0.029000000000000001)
    40 = @(1, 0)
    41 = @(0.88100000000000001, 0.45600000000000002)
    42 = @(1, 0)
    43 = @(1, 0.0070000000000000001)
    44 = @(0.68, 0.622)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

# Update the sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("I45").Select()
